$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "VALOR MORA" amount (E11): 416000 -> 468000
$ws.Range("E11").Value = 468000

# 2. Update "Cant. Periodos" count (F13): 8 -> 9
$ws.Range("F13").Value = 9

# 3. Insert a new row at 24 (pushes old rows 24+ down by one, incl. the
#    signature-block rows and their merged cells)
$ws.Rows("24").Insert()

# 4. The new row 24 becomes the new "last" period row: clone row 23's
#    content + formatting (values, borders, fills) into it.
$ws.Range("B23:J23").Copy($ws.Range("B24:J24"))

# 5. Row 23 (now a "middle" row instead of the last one) should take on
#    the non-bottom-border formatting used by the other middle rows (22).
$ws.Range("B22:J22").Copy()
$ws.Range("B23:J23").PasteSpecial(-4122)

# 6. The periods database was re-sorted to ascending order (oldest on
#    top, newest on the bottom) and the new period appended at the end.
$ws.Range("E16").Value = "2412"
$ws.Range("E17").Value = "2501"
$ws.Range("E18").Value = "2502"
$ws.Range("E19").Value = "2503"
$ws.Range("E20").Value = "2504"
$ws.Range("E21").Value = "2505"
$ws.Range("E22").Value = "2506"
$ws.Range("E23").Value = "2507"
$ws.Range("E24").Value = "2508"
